# Apply the SoCDTtiNTY -> SoCDTtiNTY-psgr / SoCDTtiNTY-frgt restructuring.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the existing data sheet to the "passenger" variant.
# ---------------------------------------------------------------------
$psgr = $wb.Worksheets.Item("SoCDTtiNTY")
$psgr.Name = "SoCDTtiNTY-psgr"

# ---------------------------------------------------------------------
# 2. Re-label the header row (A1:H1) with the new column headers and
#    the new column layout (A = row label header, B:H = vehicle techs).
# ---------------------------------------------------------------------
$psgr.Range("A1").Value = "Share that is New (dimensionless)"
$psgr.Range("B1").Value = "battery electric vehicle"
$psgr.Range("C1").Value = "natural gas vehicle"
$psgr.Range("D1").Value = "gasoline vehicle"
$psgr.Range("E1").Value = "diesel vehicle"
$psgr.Range("F1").Value = "plugin hybrid vehicle"
$psgr.Range("G1").Value = "LPG vehicle"
$psgr.Range("H1").Value = "hydrogen vehicle"

$headerRow = $psgr.Range("A1:H1")
$headerRow.WrapText = $true
$headerRow.RowHeight = 30

$psgr.Range("A1").Font.Bold = $true
$psgr.Range("B1:H1").HorizontalAlignment = -4152

# ---------------------------------------------------------------------
# 3. Column widths for the new layout.
# ---------------------------------------------------------------------
$psgr.Columns.Item(1).ColumnWidth = 18.25
$psgr.Range("B1:H1").ColumnWidth = 13.6

# ---------------------------------------------------------------------
# 4. Clear the stray formatting that used to live on B5 and fill the
#    data across the new B:H columns (same value repeated per row,
#    matching the passenger-vehicle calibration numbers).
# ---------------------------------------------------------------------
$psgr.Range("B5").ClearFormats()

$psgrValues = @{
    2 = 0.076
    3 = 0.0435
    4 = 0.0416
    5 = 0.029
    6 = 0.02982
    7 = 0.0587
}

foreach ($row in $psgrValues.Keys) {
    $value = $psgrValues[$row]
    $psgr.Range("B" + $row + ":H" + $row).Value = $value
}

# ---------------------------------------------------------------------
# 5. Duplicate the fully-formatted passenger sheet to create the
#    freight sheet, then overwrite its data values.
# ---------------------------------------------------------------------
$psgr.Copy([System.Reflection.Missing]::Value, $psgr)
$frgt = $wb.Worksheets.Item($psgr.Index + 1)
$frgt.Name = "SoCDTtiNTY-frgt"

$frgtValues = @{
    2 = 0.07
    3 = 0.035
    4 = 0.042
    5 = 0.029
    6 = 0.0303
    7 = 0
}

foreach ($row in $frgtValues.Keys) {
    $value = $frgtValues[$row]
    $frgt.Range("B" + $row + ":H" + $row).Value = $value
}
